# Generate Report for Handback
# Fills in the "Latest Target File", "Latest Handback File" and
# "Latest Handback DateTime" columns for each locale sheet, and flips the
# Status column from "Ready for handoff" to "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$mdBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5f6517a0e9a649fbacd72a3585d3055a7739fd96/e2e/"
$file1 = "436c0d43-49b8-48e9-b962-c3ad530268c7.md"
$file2 = "df4c6466-47dd-4ebd-8626-78028a3eb273.md"

# ---------------------------------------------------------------------
# Overview sheet: Status shows up in the zh-cn / de-de summary columns.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusNew
$wsOverview.Range("F2").Value = $statusNew
$wsOverview.Range("E3").Value = $statusNew
$wsOverview.Range("F3").Value = $statusNew
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusNew
$wsZh.Range("C3").Value = $statusNew

$wsZh.Range("J2").Value = "436c0d43-49b8-48e9-b962-c3ad530268c7.559b2dc02666d81a511cabbf18537dc5e780a6d9.zh-cn.xlf"
$wsZh.Range("J3").Value = "df4c6466-47dd-4ebd-8626-78028a3eb273.715ac8436e5aee92cc67e873fc94adb7dfb62cec.zh-cn.xlf"

$wsZh.Range("K2").Value = "2016-08-24 14:32:28"
$wsZh.Range("K3").Value = "2016-08-24 14:32:28"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $mdBase + $file1, "", "", $file1)
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdBase + $file1, "", "", $file1)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $mdBase + $file2, "", "", $file2)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $mdBase + $file2, "", "", $file2)

$wsZh.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZh.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZh.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusNew
$wsDe.Range("C3").Value = $statusNew

$wsDe.Range("J2").Value = "436c0d43-49b8-48e9-b962-c3ad530268c7.559b2dc02666d81a511cabbf18537dc5e780a6d9.de-de.xlf"
$wsDe.Range("J3").Value = "df4c6466-47dd-4ebd-8626-78028a3eb273.715ac8436e5aee92cc67e873fc94adb7dfb62cec.de-de.xlf"

$wsDe.Range("K2").Value = "2016-08-24 14:32:35"
$wsDe.Range("K3").Value = "2016-08-24 14:32:35"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $mdBase + $file1, "", "", $file1)
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdBase + $file1, "", "", $file1)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $mdBase + $file2, "", "", $file2)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $mdBase + $file2, "", "", $file2)

$wsDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDe.Columns.Item(10).ColumnWidth = 39.166666666666664
